# Fixing the big mistake: update computed average/statistics values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: mean
$ws.Range("B3").Value = 44694.22931324951
$ws.Range("D3").Value = 2566.918741906393
$ws.Range("E3").Value = 7690.107310439496

# Row 4: std
$ws.Range("B4").Value = 19209.18146037113
$ws.Range("D4").Value = 1337.163928315268
$ws.Range("E4").Value = 3445.543304400405

# Row 5: min
$ws.Range("B5").Value = 16160.40429589045
$ws.Range("D5").Value = 478.4026849315056
$ws.Range("E5").Value = 1992.98733150685

# Row 6: 25%
$ws.Range("B6").Value = 28966.08085068493
$ws.Range("D6").Value = 1303.408328082193
$ws.Range("E6").Value = 4040.291458219175

# Row 7: 50%
$ws.Range("B7").Value = 41860.85765479452
$ws.Range("D7").Value = 2217.218943835618
$ws.Range("E7").Value = 7600.181989041095

# Row 8: 75%
$ws.Range("B8").Value = 58712.55744657514
$ws.Range("D8").Value = 3737.860337671233
$ws.Range("E8").Value = 11401.70497465753

# Row 9: max
$ws.Range("B9").Value = 80499.78310136932
$ws.Range("D9").Value = 4796.593706849323
$ws.Range("E9").Value = 13288.28970410958

# Row 10: Total sum
$ws.Range("F10").Value = 64359690.21107931

# Row 11: Residential % energy sector
$ws.Range("G11").Value = 0.7705067027678848

# Row 12: Community total + % energy sector
$ws.Range("F12").Value = 3696362.988345208
$ws.Range("G12").Value = 0.05743288968952947

# Row 13: IGA total + % energy sector
$ws.Range("F13").Value = 11073754.52703287
$ws.Range("G13").Value = 0.1720604075425858
